# clarify status report for unit tests
#
# The "B" column cells that previously cached a QuantLib-XLL-computed
# "cvNN#0000" object-handle string now fail with #NUM! (the add-in call
# itself errors out). The "D" column's PASS/FAIL formula is rewritten so
# that an error in either operand is reported as "ERROR" instead of being
# silently compared (and, for B8/B9 which already errored, surfaced as
# "ERROR" rather than propagating the raw #N/A / #NUM! error).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B column: these calls now error out with #NUM! ------------------
$ws.Range("B3").Value  = "#NUM!"
$ws.Range("B4").Value  = "#NUM!"
$ws.Range("B5").Value  = "#NUM!"
$ws.Range("B6").Value  = "#NUM!"
$ws.Range("B7").Value  = "#NUM!"
$ws.Range("B10").Value = "#NUM!"
$ws.Range("B11").Value = "#NUM!"

# --- D column: report errors explicitly instead of comparing them ----
# D3 is entered on its own (stays an individual, non-shared formula);
# D4:D15 are filled down together as one shared-formula block.
$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'
$ws.Range("D4:D15").Formula = '=IF(ISERROR(B4),"ERROR",IF(ISERROR(C4),"FAIL",IF(B4=C4,"PASS","FAIL")))'
